# Scheduled-runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, Leve price/profit columns) across the per-job
# "Zeromus_Profits" sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR). A few rows
# whose listings disappeared from the market board lose their profit cell
# entirely (ClearContents), and one row (GSM!N22) gains a freshly computed
# HQ-profit cell that didn't exist before.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 220
$ws.Range("I18").Value = 220
$ws.Range("K18").Value = 220
$ws.Range("M18").Value = 64
$ws.Range("H40").Value = 3252.8333
$ws.Range("I40").Value = 7760
$ws.Range("K40").Value = 7760
$ws.Range("M40").Value = -7585
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
$ws.Range("H141").Value = 3137.325
$ws.Range("I141").Value = 1803.3214
$ws.Range("J141").Value = 6250
$ws.Range("K141").Value = 5409.9642
$ws.Range("L141").Value = 18750
$ws.Range("M141").Value = -229.9642000000003
$ws.Range("N141").Value = -29110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 40869.46
$ws.Range("I32").Value = 33647.12
$ws.Range("J32").Value = 47879.383
$ws.Range("K32").Value = 33647.12
$ws.Range("L32").Value = 47879.383
$ws.Range("M32").Value = -33360.12
$ws.Range("N32").Value = -48453.383
$ws.Range("H102").Value = 1602.4445
$ws.Range("I102").Value = 1250
$ws.Range("J102").Value = 1703.1428
$ws.Range("K102").Value = 1250
$ws.Range("L102").Value = 1703.1428
$ws.Range("M102").Value = 372
$ws.Range("N102").Value = -4947.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 62164.445
$ws.Range("J4").Value = 7435
$ws.Range("L4").Value = 7435
$ws.Range("N4").Value = -7659
$ws.Range("H21").Value = 5738.3335
$ws.Range("J21").Value = 7607.5
$ws.Range("L21").Value = 7607.5
$ws.Range("N21").Value = -8077.5
$ws.Range("H31").Value = 2152.5356
$ws.Range("I31").Value = 1237.9474
$ws.Range("J31").Value = 4083.3333
$ws.Range("K31").Value = 1237.9474
$ws.Range("L31").Value = 4083.3333
$ws.Range("M31").Value = -942.9474
$ws.Range("N31").Value = -4673.3333
$ws.Range("H34").Value = 2152.5356
$ws.Range("I34").Value = 1237.9474
$ws.Range("J34").Value = 4083.3333
$ws.Range("K34").Value = 1237.9474
$ws.Range("L34").Value = 4083.3333
$ws.Range("M34").Value = -1035.9474
$ws.Range("N34").Value = -4487.3333
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H41").Value = 2029.5
$ws.Range("I41").Value = 2029.5
$ws.Range("K41").Value = 2029.5
$ws.Range("M41").Value = -1601.5
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H50").Value = 18900
$ws.Range("J50").Value = 18900
$ws.Range("L50").Value = 18900
$ws.Range("N50").Value = -20150
$ws.Range("H51").Value = 13233.333
$ws.Range("J51").Value = 13233.333
$ws.Range("L51").Value = 13233.333
$ws.Range("N51").Value = -14705.333
$ws.Range("H59").Value = 19680
$ws.Range("J59").Value = 19680
$ws.Range("L59").Value = 19680
$ws.Range("N59").Value = -21970
$ws.Range("H61").Value = 13233.333
$ws.Range("J61").Value = 13233.333
$ws.Range("L61").Value = 13233.333
$ws.Range("N61").Value = -13929.333
$ws.Range("H68").Value = 19500
$ws.Range("J68").Value = 19500
$ws.Range("L68").Value = 19500
$ws.Range("N68").Value = -20998
$ws.Range("H71").Value = 19500
$ws.Range("J71").Value = 19500
$ws.Range("L71").Value = 58500
$ws.Range("N71").Value = -65988
$ws.Range("H132").Value = 1538.6136
$ws.Range("I132").Value = 969.1
$ws.Range("J132").Value = 2759
$ws.Range("K132").Value = 2907.3
$ws.Range("L132").Value = 8277
$ws.Range("M132").Value = -377.3000000000002
$ws.Range("N132").Value = -13337

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H19").Value = 963.1579
$ws.Range("I19").Value = 963.1579
$ws.Range("K19").Value = 2889.4737
$ws.Range("M19").Value = -2715.4737
$ws.Range("H104").Value = 1993.2727
$ws.Range("I104").Value = 1426
$ws.Range("J104").Value = 2050
$ws.Range("K104").Value = 4278
$ws.Range("L104").Value = 6150
$ws.Range("M104").Value = -1657
$ws.Range("N104").Value = -11392
$ws.Range("H129").Value = 16667680
$ws.Range("I129").Value = 402.5
$ws.Range("J129").Value = 27779198
$ws.Range("K129").Value = 1207.5
$ws.Range("L129").Value = 83337594
$ws.Range("M129").Value = 3792.5
$ws.Range("N129").Value = -83347594
$ws.Range("H132").Value = 880.9677
$ws.Range("I132").Value = 565.94116
$ws.Range("J132").Value = 1263.5
$ws.Range("K132").Value = 5093.47044
$ws.Range("L132").Value = 11371.5
$ws.Range("M132").Value = -2563.47044
$ws.Range("N132").Value = -16431.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 7997.273
$ws.Range("I5").Value = 2660
$ws.Range("J5").Value = 9998.75
$ws.Range("K5").Value = 2660
$ws.Range("L5").Value = 9998.75
$ws.Range("M5").Value = -2548
$ws.Range("N5").Value = -10222.75
$ws.Range("H22").Value = 1625
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1833.3334
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1833.3334
$ws.Range("M22").Value = -471
$ws.Range("N22").Value = -2891.3334
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H133").Value = 33000
$ws.Range("J133").Value = 33000
$ws.Range("L133").Value = 33000
$ws.Range("N133").Value = -43120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1504687.5
$ws.Range("J2").Value = 1005357.1
$ws.Range("L2").Value = 1005357.1
$ws.Range("N2").Value = -1005581.1
$ws.Range("H46").Value = 2266.6667
$ws.Range("J46").Value = 800
$ws.Range("L46").Value = 800
$ws.Range("N46").Value = -1176
$ws.Range("H68").Value = 30760032
$ws.Range("I68").Value = 42292732
$ws.Range("J68").Value = 6166.6665
$ws.Range("K68").Value = 42292732
$ws.Range("L68").Value = 6166.6665
$ws.Range("M68").Value = -42291983
$ws.Range("N68").Value = -7664.6665
$ws.Range("H71").Value = 30760032
$ws.Range("I71").Value = 42292732
$ws.Range("J71").Value = 6166.6665
$ws.Range("K71").Value = 211463660
$ws.Range("L71").Value = 30833.3325
$ws.Range("M71").Value = -211459916
$ws.Range("N71").Value = -38321.3325

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 8939
$ws.Range("I2").Value = 5500
$ws.Range("J2").Value = 9430.286
$ws.Range("L2").Value = 9430.286
$ws.Range("N2").Value = -9654.286
$ws.Range("H23").Value = 1036.6666
$ws.Range("I23").Value = 1036.6666
$ws.Range("K23").Value = 1036.6666
$ws.Range("M23").Value = -807.6666
